# Updated cryptos list on Tue Oct 24 23:49:18 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeRef, $text) {
    $rng = $ws.Range($rangeRef)
    # Force text interpretation so values like "225.85" or "0.625" are not
    # auto-converted to numbers, then restore the original (default) style
    # so no new cell formatting is introduced.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "33.866.76"
Set-TextValue "E2" "  +3.20%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.786.76"
Set-TextValue "E3" "  +1.57%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.37%  "

# Row 5 - BNB
Set-TextValue "D5" "225.85"
Set-TextValue "E5" "  -0.54%  "

# Row 6 - XRP
Set-TextValue "E6" "  +2.61%  "

# Row 7 - USDC
Set-TextValue "E7" "  +0.39%  "

# Row 8 - Solana
Set-TextValue "D8" "30.14"
Set-TextValue "E8" "  -5.00%  "

# Row 9 - OKB
Set-TextValue "D9" "46.60"
Set-TextValue "E9" "  +2.98%  "

# Row 10 - Cardano
Set-TextValue "E10" "  +0.45%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.0666"
Set-TextValue "E11" "  -0.26%  "

# Row 12 - TRON
Set-TextValue "E12" "  +0.83%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "2.044.22"
Set-TextValue "E13" "  +1.71%  "

# Row 14 - WrappedEther
Set-TextValue "D14" "1.783.66"
Set-TextValue "E14" "  +1.54%  "

# Row 15 - was Polygon, now Chainlink
Set-TextValue "B15" "Chainlink"
Set-TextValue "C15" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D15" "10.46"
Set-TextValue "E15" "  +0.17%  "

# Row 16 - was Chainlink, now Polygon
Set-TextValue "B16" "Polygon"
Set-TextValue "C16" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D16" "0.625"
Set-TextValue "E16" "  -0.93%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "33.888.33"
Set-TextValue "E17" "  +3.43%  "

# Row 18 - Polkadot
Set-TextValue "E18" "  -2.50%  "

# Row 19 - Litecoin
Set-TextValue "D19" "69.02"
Set-TextValue "E19" "  +0.40%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "251.80"
Set-TextValue "E20" "  -2.48%  "

# Row 21 - ShibaInu
Set-TextValue "E21" "  -0.64%  "

# Row 22 - Dai
Set-TextValue "E22" "  +0.50%  "

# Row 23 - Avalanche
Set-TextValue "E23" "  -1.38%  "

# Row 24 - Uniswap
Set-TextValue "D24" "4.22"
Set-TextValue "E24" "  -3.15%  "

# Row 25 - Toncoin
Set-TextValue "E25" "  -1.81%  "

# Row 26 - Monero
Set-TextValue "D26" "158.17"
Set-TextValue "E26" "  -0.88%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "16.45"
Set-TextValue "E27" "  -0.34%  "

# Row 28 - Stellar
Set-TextValue "E28" "  -1.02%  "

# Row 29 - Cosmos
Set-TextValue "D29" "6.98"
Set-TextValue "E29" "  +0.14%  "

# Row 30 - BinanceUSD
Set-TextValue "E30" "  +0.29%  "

# Row 31 - Filecoin
Set-TextValue "E31" "  -0.72%  "

# Row 32 - Hedera
Set-TextValue "E32" "  -0.33%  "

# Row 33 - PancakeSwap
Set-TextValue "E33" "  +1.52%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue "D34" "3.59"
Set-TextValue "E34" "  +2.92%  "

# Row 35 - LidoDAOToken
Set-TextValue "E35" "  +3.58%  "

# Row 36 - Maker
Set-TextValue "D36" "1.501.13"
Set-TextValue "E36" "  -3.61%  "

# Row 37 - TrustWalletToken
Set-TextValue "D37" "1.06"
Set-TextValue "E37" "  +1.74%  "

# Row 38 - ImmutableX
Set-TextValue "D38" "0.633"
Set-TextValue "E38" "  +0.10%  "

# Row 39 - was VeChain, now Aave
Set-TextValue "B39" "Aave"
Set-TextValue "C39" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D39" "83.33"
Set-TextValue "E39" "  -1.30%  "

# Row 40 - was Aave, now VeChain
Set-TextValue "B40" "VeChain"
Set-TextValue "C40" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D40" "0.0185"
Set-TextValue "E40" "  -0.11%  "

# Row 41 - HuobiToken
Set-TextValue "D41" "2.35"
Set-TextValue "E41" "  +1.87%  "

# Row 42 - MXToken
Set-TextValue "E42" "  -0.88%  "

# Row 43 - ARBITRUM
Set-TextValue "E43" "  +2.82%  "

# Row 44 - Kaspa
Set-TextValue "E44" "  +0.93%  "

# Row 45 - RenderToken
Set-TextValue "E45" "  -1.98%  "

# Row 46 - WEMIXToken
Set-TextValue "E46" "  +2.08%  "

# Row 47 - RocketPoolETH
Set-TextValue "D47" "1.935.96"
Set-TextValue "E47" "  +1.42%  "

# Row 48 - FraxShare
Set-TextValue "D48" "5.71"
Set-TextValue "E48" "  +0.63%  "

# Row 49 - PaxDollar
Set-TextValue "E49" "  +0.30%  "

# Row 50 - InjectiveProtocol
Set-TextValue "D50" "11.81"
Set-TextValue "E50" "  +8.63%  "

# Row 51 - BitcoinSV
Set-TextValue "D51" "51.25"
Set-TextValue "E51" "  -5.05%  "
